# Slide 3 ("TEAM MEMBERS" banner): widen the orange rounded-rectangle
# background and the "TEAM MEMBERS" text box that sits on top of it, and
# let the text box wrap instead of forcing a single line ("wrap=none" ->
# "wrap=square").
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

$roundRect = $s.Shapes.Item(11)   # "圆角矩形 14" (rounded-rectangle background)
$teamLabel = $s.Shapes.Item(12)   # "矩形 15" (the "TEAM MEMBERS" text box)

# Rounded rectangle only grows in width; position/height stay put.
# NOTE: PowerPoint's Shape.Width/Left/Top are single-precision (float32)
# COM properties, so a literal 374.4 would truncate to 4754879 EMU instead
# of the intended 4754880 EMU. Nudge to the nearest representable float32
# that lands on the exact EMU target.
$roundRect.Width = 374.4000244140625   # 4754880 EMU (was 3870960 EMU / 304.8pt)

# Text box moves right/down slightly and grows in width to match.
$teamLabel.Left = 109.4    # 1389380 EMU (was 1022073 EMU)
$teamLabel.Top = 68.9      # 875030 EMU (was 874737 EMU)
$teamLabel.Width = 331.85  # 4214495 EMU (was 3261360 EMU)

# Allow the text to wrap within the (now wider) box instead of never wrapping.
$teamLabel.TextFrame.WordWrap = -1   # msoTrue -> wrap="square"
